# The sheet holds one weekly price-report row per (date, quality) pair.
# This edit inserts one new weekly data point (a new "Primera" quality
# observation dated 45194) right after the existing row 1113, pushing every
# subsequent row down by one (old row 1114 -> new row 1115, ..., old row
# 1233 -> new row 1234). The sheet's used range grows from A1:R1233 to
# A1:R1234 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 1114, shifting rows 1114:1233 down to 1115:1234.
$ws.Rows(1114).Insert()

# Populate the newly inserted row 1114 with the new observation.
$ws.Range("A1114").Value = 3
$ws.Range("B1114").Value = "Femacal de La Calera"
$ws.Range("C1114").Value = "Coquimbo"
$ws.Range("D1114").Value = 45194
$ws.Range("E1114").Value = 5
$ws.Range("F1114").Value = 100114014
$ws.Range("G1114").Value = "Betarraga"
$ws.Range("H1114").Value = "Sin especificar"
$ws.Range("I1114").Value = "Primera"
$ws.Range("J1114").Value = 3400
$ws.Range("K1114").Value = 500
$ws.Range("L1114").Value = 550
$ws.Range("M1114").Value = 524
$ws.Range("N1114").Value = "`$/paquete 4 unidades"
$ws.Range("O1114").Value = "Provincia de Quillota"
$ws.Range("P1114").Value = 131
$ws.Range("Q1114").Value = 4
$ws.Range("R1114").Value = "Hortaliza"
